# Ajuste de la Queratina / % a Marinela / Menor % a los profesionales de Tocador
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert columns ---
# Step 1: insert 3 new columns before F (old F,G,H "Valor_producto","Part_profesional","Revisar"
#         shift right to I,J,K)
$ws.Range("F1:H1").EntireColumn.Insert() | Out-Null

# Step 2: insert 1 new column before J (currently "Part_profesional"), shifting it to K
#         and "Revisar" from K to L
$ws.Range("J1").EntireColumn.Insert() | Out-Null

# --- New header row values (Insert carries the existing header style s="1" along) ---
$ws.Range("F1").Value = "Porc_trans"
$ws.Range("G1").Value = "Cost_trans"
$ws.Range("H1").Value = "Porc_producto"
$ws.Range("J1").Value = "Valor_Neto"

# --- Row 2 data ---
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.1166
$ws.Range("J2").Value = 30919

# --- Row 3 data ---
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.2720666666666667
$ws.Range("J3").Value = 10919
